$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$NL = [char]10

# --- Row 3: fix typo "invalue" -> "invalid" in the Solution/Workaround cell, and shrink row height ---
$ws.Range("B3").Value = "Generally, try to find and tweak the string that may be causing the issue--usually it's because GIS is trying to confirm if it's a data source, and if it finds it is invalid it will through the error." + $NL + "Workaround: make one offending url string a 1-item list, then pluck it back out of the list once it's being used--that way, ArcGIS thinks it's a list an doesn't scrutinize it like a string."
$ws.Rows(3).RowHeight = 105

# --- Row 4: add a new "ESRI Follow-up?" hyperlink cell (A4:C4 are unchanged) ---
$ws.Range("D4").Value = "Geonet thread DC posted 3/31/2022"
$ws.Hyperlinks.Add($ws.Range("D4"), "https://community.esri.com/t5/arcgis-pro-questions/geonet-thread-dc-posted-3-31-2022/td-p/000358", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value) | Out-Null
$ws.Range("D4").WrapText = $true

# --- Row 5: brand new troubleshooting-log entry ---
$ws.Range("A5").Value = "000358: Invalid expression when running select-by-attributes command in arcpy"
$ws.Range("B5").Value = "If it happens when running in Arc Pro, make sure that none of the columns are hidden for the layer that the selection is happening on."
$ws.Range("C5").Value = "NA"
$ws.Range("D5").Value = "NA"
$ws.Rows(5).RowHeight = 45

# --- Selection moves to B5 after the edit ---
$ws.Range("B5").Select() | Out-Null
